$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (unchanged values, just confirm them)
$ws.Range("A1").Value = "tech_id"
$ws.Range("B1").Value = "tech_name"
$ws.Range("C1").Value = "username"
$ws.Range("D1").Value = "pwd"

# Row 2 - Berlin
$ws.Range("A2").Value = 597041
$ws.Range("B2").Value = "Berlin"
$ws.Range("C2").Value = "berlin12"
$ws.Range("D2").Style = "Normal"
$ws.Range("D2").Value = "123456789a"

# Row 3 - Helsinki
$ws.Range("A3").Value = 748836
$ws.Range("B3").Value = "Helsinki"
$ws.Range("C3").Value = "helsinki"
$ws.Range("D3").Style = "Normal"
$ws.Range("D3").Value = "123456789a"

# Row 4 - Cairo
$ws.Range("A4").Value = 330228
$ws.Range("B4").Value = "Cairo"
$ws.Range("C4").Value = "cairo123"
$ws.Range("D4").Style = "Normal"
$ws.Range("D4").Value = "123456789a"

# Row 5 - Paris (new row)
$ws.Range("A5").Value = 151077
$ws.Range("B5").Value = "Paris"
$ws.Range("C5").Value = "paris123"
$ws.Range("D5").Value = "123456789a"

$ws.Range("B4").Select()
